$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.722.85"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.601.02"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.19"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0620"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "1.825.28"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.601.30"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.19"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "26.690.54"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.84"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.02"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.10"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0512"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "1.294.98"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.607"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.18"
$ws.Range("E38").Value = "  +17.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0170"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.823"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.779"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.28"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "1.737.92"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.84"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.41"
$ws.Range("E51").Value = "  +0.15%  "
